$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 2
$ws.Range("P2").Value = 16.25
$ws.Range("Q2").Value = 8.049847157796222
$ws.Range("R2").Value = 7.492214147796222

# Delete column Z entirely (header Z1 and value Z2), shrinking dimension to A1:Y2
$ws.Range("Z1:Z2").Delete()
